# Apply the latest cryptos snapshot values scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.446.96"
$ws.Range("D3").Value = "'2.654.01"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'602.24"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'157.33"
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("E9").Value = "  +6.60%  "
$ws.Range("D10").Value = "'0.403"
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("D11").Value = "'5.83"
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "'29.27"
$ws.Range("E13").Value = "  +4.72%  "
$ws.Range("D14").Value = "'3.128.74"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").Value = "'65.257.82"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "'0.0000174"
$ws.Range("E16").Value = "  +12.29%  "
$ws.Range("D17").Value = "'2.635.25"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'12.62"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "'4.84"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").Value = "'355.04"
$ws.Range("D21").Value = "'7.28"
$ws.Range("E21").Value = "  +5.62%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "'68.21"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'9.61"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "'8.38"
$ws.Range("E27").Value = "  +4.19%  "
$ws.Range("D28").Value = "'0.166"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").Value = "'541.60"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").Value = "'0.994"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "'0.0₃0926"
$ws.Range("E31").Value = "  +8.46%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "'1.84"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("D34").Value = "'5.82"
$ws.Range("E34").Value = "  +9.36%  "
$ws.Range("D35").Value = "'6.47"
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("D36").Value = "'0.429"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("E37").Value = "  +5.84%  "
$ws.Range("D38").Value = "'165.44"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "'20.26"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'168.65"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("E43").Value = "  +4.97%  "
$ws.Range("D44").Value = "'4.12"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0610"
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'23.43"
$ws.Range("E46").Value = "  +5.77%  "
$ws.Range("E47").Value = "  +10.88%  "
$ws.Range("D48").Value = "'0.649"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "'0.0253"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "'19.56"
$ws.Range("E51").Value = "  +1.14%  "
